# Update a handful of the (randomly-generated) sample readings on the
# "Gearbox Tests" sheet. This accompanies an unrelated repo-wide
# style/tooling cleanup (black/flake8/isort, venv bump) -- the sheet
# values below are simply refreshed fixture data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gearbox Tests")

$ws.Range("C7").Value = 54.455593060061851
$ws.Range("C8").Value = 90
$ws.Range("C11").Value = 899.99999999999909
